$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1822660098522167
$ws.Range("C2").Value = 0.5467980295566502
$ws.Range("J2").Value = 0.009852216748768473
$ws.Range("P2").Value = 0.1625615763546798
$ws.Range("S2").Value = 0.09852216748768473
$ws.Range("C3").Value = 0.02631578947368421
$ws.Range("J3").Value = 0.04385964912280702
$ws.Range("P3").Value = 0.7192982456140351
$ws.Range("S3").Value = 0.2105263157894737
$ws.Range("J4").Value = 0.04347826086956522
$ws.Range("P4").Value = 0.6739130434782609
$ws.Range("S4").Value = 0.2826086956521739
$ws.Range("B6").Value = 0.05164319248826291
$ws.Range("D6").Value = 0.01408450704225352
$ws.Range("F6").Value = 0.03755868544600939
$ws.Range("J6").Value = 0.2206572769953052
$ws.Range("O6").Value = 0.004694835680751174
$ws.Range("Q6").Value = 0.1267605633802817
$ws.Range("R6").Value = 0.07042253521126761
$ws.Range("S6").Value = 0.4741784037558686
$ws.Range("B7").Value = 0.05084745762711865
$ws.Range("D7").Value = 0.03954802259887006
$ws.Range("E7").Value = 0.005649717514124294
$ws.Range("F7").Value = 0.05649717514124294
$ws.Range("J7").Value = 0.1186440677966102
$ws.Range("O7").Value = 0.01694915254237288
$ws.Range("Q7").Value = 0.1694915254237288
$ws.Range("R7").Value = 0.1299435028248588
$ws.Range("S7").Value = 0.4124293785310734
$ws.Range("B8").Value = 0.06088560885608856
$ws.Range("D8").Value = 0.01476014760147601
$ws.Range("F8").Value = 0.07195571955719557
$ws.Range("J8").Value = 0.08856088560885608
$ws.Range("O8").Value = 0.02583025830258303
$ws.Range("Q8").Value = 0.1752767527675277
$ws.Range("R8").Value = 0.0940959409594096
$ws.Range("S8").Value = 0.4686346863468634
$ws.Range("B9").Value = 0.05853658536585366
$ws.Range("D9").Value = 0.01463414634146342
$ws.Range("F9").Value = 0.06341463414634146
$ws.Range("J9").Value = 0.08292682926829269
$ws.Range("O9").Value = 0.01463414634146342
$ws.Range("Q9").Value = 0.1414634146341463
$ws.Range("R9").Value = 0.09268292682926829
$ws.Range("S9").Value = 0.5317073170731708
$ws.Range("B10").Value = 0.07604251839738348
$ws.Range("D10").Value = 0.02207686017988553
$ws.Range("F10").Value = 0.06295993458708095
$ws.Range("J10").Value = 0.1234668847097302
$ws.Range("O10").Value = 0.01062959934587081
$ws.Range("Q10").Value = 0.2044153720359771
$ws.Range("R10").Value = 0.07849550286181521
$ws.Range("S10").Value = 0.4219133278822568
$ws.Range("G11").Value = 0.1660079051383399
$ws.Range("J11").Value = 0.05928853754940711
$ws.Range("K11").Value = 0.1976284584980237
$ws.Range("L11").Value = 0.5652173913043478
$ws.Range("S11").Value = 0.01185770750988142
$ws.Range("G12").Value = 0.7549668874172185
$ws.Range("J12").Value = 0.1589403973509934
$ws.Range("K12").Value = 0.01324503311258278
$ws.Range("L12").Value = 0.01986754966887417
$ws.Range("S12").Value = 0.05298013245033113
$ws.Range("G13").Value = 0.7297297297297297
$ws.Range("J13").Value = 0.2432432432432433
$ws.Range("S13").Value = 0.02702702702702703
$ws.Range("F15").Value = 0.02185792349726776
$ws.Range("H15").Value = 0.180327868852459
$ws.Range("I15").Value = 0.1256830601092896
$ws.Range("J15").Value = 0.3060109289617486
$ws.Range("K15").Value = 0.04371584699453552
$ws.Range("M15").Value = 0.01639344262295082
$ws.Range("N15").Value = 0.00546448087431694
$ws.Range("O15").Value = 0.0273224043715847
$ws.Range("S15").Value = 0.273224043715847
$ws.Range("F16").Value = 0.007092198581560284
$ws.Range("H16").Value = 0.2340425531914894
$ws.Range("I16").Value = 0.148936170212766
$ws.Range("J16").Value = 0.3404255319148936
$ws.Range("K16").Value = 0.09929078014184398
$ws.Range("N16").Value = 0.007092198581560284
$ws.Range("O16").Value = 0.05673758865248227
$ws.Range("S16").Value = 0.1063829787234043
$ws.Range("F17").Value = 0.02612826603325416
$ws.Range("H17").Value = 0.2280285035629454
$ws.Range("I17").Value = 0.1211401425178147
$ws.Range("J17").Value = 0.3942992874109263
$ws.Range("K17").Value = 0.0831353919239905
$ws.Range("M17").Value = 0.01425178147268409
$ws.Range("O17").Value = 0.0498812351543943
$ws.Range("S17").Value = 0.0831353919239905
$ws.Range("F18").Value = 0.009950248756218905
$ws.Range("H18").Value = 0.1940298507462687
$ws.Range("I18").Value = 0.0845771144278607
$ws.Range("J18").Value = 0.4278606965174129
$ws.Range("K18").Value = 0.07960199004975124
$ws.Range("M18").Value = 0.01492537313432836
$ws.Range("O18").Value = 0.04477611940298507
$ws.Range("S18").Value = 0.1442786069651741
$ws.Range("F19").Value = 0.01138790035587189
$ws.Range("H19").Value = 0.2476868327402135
$ws.Range("I19").Value = 0.06548042704626335
$ws.Range("J19").Value = 0.3886120996441281
$ws.Range("K19").Value = 0.0903914590747331
$ws.Range("M19").Value = 0.01779359430604982
$ws.Range("O19").Value = 0.05907473309608541
$ws.Range("S19").Value = 0.1195729537366548
